# Fix total marks error on the "quiz" marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (row 11): Right count 5 -> 4, Wrong count -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Total row (row 12): Right total 70 -> 56, and the "x / y" summary text
$ws.Range("B12").Value = 56
$ws.Range("E12").Value = "56 / 112"
